# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition data) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 147
$ws1.Range("F3").Value = 462
$ws1.Range("F4").Value = 7
$ws1.Range("F7").Value = 26
$ws1.Range("F8").Value = 10
$ws1.Range("F9").Value = 167
$ws1.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202408/V8EaWtul1724135798600.jpeg"

# --- Sheet "演出" (performance data) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 81
$ws2.Range("F3").Value = 33
$ws2.Range("F4").Value = 2

# --- Sheet "全部类型" (combined / all types data) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 147
$ws4.Range("F3").Value = 81
$ws4.Range("F4").Value = 462
$ws4.Range("F5").Value = 7
$ws4.Range("F8").Value = 26
$ws4.Range("F9").Value = 10
$ws4.Range("F10").Value = 167
$ws4.Range("F11").Value = 33
$ws4.Range("F12").Value = 2
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202408/V8EaWtul1724135798600.jpeg"
